$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H44").Value = 6000
$ws.Range("I44").Value = 6000
$ws.Range("K44").Value = 6000
$ws.Range("M44").Value = -5538
$ws.Range("H69").Value = 2882.9546
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 2882.9546
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 8648.863799999999
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -10396.8638
$ws.Range("H72").Value = 2882.9546
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 2882.9546
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 25946.5914
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -34682.5914
$ws.Range("H129").Value = 11178.048
$ws.Range("I129").Value = 299.375
$ws.Range("J129").Value = 17872.615
$ws.Range("K129").Value = 898.125
$ws.Range("L129").Value = 53617.845
$ws.Range("M129").Value = 4101.875
$ws.Range("N129").Value = -63617.845
$ws.Range("H132").Value = 4152.2383
$ws.Range("I132").Value = 4386.2666
$ws.Range("J132").Value = 3567.1667
$ws.Range("K132").Value = 13158.7998
$ws.Range("L132").Value = 10701.5001
$ws.Range("M132").Value = -10628.7998
$ws.Range("N132").Value = -15761.5001
$ws.Range("H137").Value = 1490.4412
$ws.Range("I137").Value = 1262.091
$ws.Range("J137").Value = 1909.0834
$ws.Range("K137").Value = 3786.273
$ws.Range("L137").Value = 5727.2502
$ws.Range("M137").Value = -1236.273
$ws.Range("N137").Value = -10827.2502

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H21").Value = 2017
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H45").Value = 1910.7778
$ws.Range("I45").Value = 1944.2632
$ws.Range("K45").Value = 1944.2632
$ws.Range("M45").Value = -1567.2632
$ws.Range("H61").Value = 1762
$ws.Range("I61").Value = 1289.5
$ws.Range("J61").Value = 2707
$ws.Range("K61").Value = 1289.5
$ws.Range("L61").Value = 2707
$ws.Range("M61").Value = -1077.5
$ws.Range("N61").Value = -3131
$ws.Range("H74").Value = 4502.4736
$ws.Range("I74").Value = 5150.6
$ws.Range("J74").Value = 3782.3333
$ws.Range("K74").Value = 5150.6
$ws.Range("L74").Value = 3782.3333
$ws.Range("M74").Value = -4276.6
$ws.Range("N74").Value = -5530.3333
$ws.Range("H77").Value = 4502.4736
$ws.Range("I77").Value = 5150.6
$ws.Range("J77").Value = 3782.3333
$ws.Range("K77").Value = 25753
$ws.Range("L77").Value = 18911.6665
$ws.Range("M77").Value = -21385
$ws.Range("N77").Value = -27647.6665
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H110").Value = 1439.25
$ws.Range("I110").Value = 1267.2941
$ws.Range("J110").Value = 1856.8572
$ws.Range("K110").Value = 1267.2941
$ws.Range("L110").Value = 1856.8572
$ws.Range("M110").Value = 777.7058999999999
$ws.Range("N110").Value = -5946.8572
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 2747.3157
$ws.Range("I132").Value = 2499.9
$ws.Range("J132").Value = 3022.2222
$ws.Range("K132").Value = 7499.700000000001
$ws.Range("L132").Value = 9066.6666
$ws.Range("M132").Value = -4969.700000000001
$ws.Range("N132").Value = -14126.6666
$ws.Range("H135").Value = 30000000
$ws.Range("J135").Value = 30000000
$ws.Range("L135").Value = 30000000
$ws.Range("N135").Value = -30010140
$ws.Range("H136").Value = 1762
$ws.Range("I136").Value = 1289.5
$ws.Range("J136").Value = 2707
$ws.Range("K136").Value = 3868.5
$ws.Range("L136").Value = 8121
$ws.Range("M136").Value = -1318.5
$ws.Range("N136").Value = -13221

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H134").Value = 1828.5238
$ws.Range("I134").Value = 1231.909
$ws.Range("J134").Value = 2484.8
$ws.Range("K134").Value = 3695.727
$ws.Range("L134").Value = 7454.400000000001
$ws.Range("M134").Value = -1160.727
$ws.Range("N134").Value = -12524.4

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H16").Value = 1843.1666
$ws.Range("I16").Value = 2566.3333
$ws.Range("J16").Value = 1120
$ws.Range("K16").Value = 2566.3333
$ws.Range("L16").Value = 1120
$ws.Range("M16").Value = -2279.3333
$ws.Range("N16").Value = -1694
$ws.Range("H31").Value = 3450751.5
$ws.Range("I31").Value = 11112581
$ws.Range("J31").Value = 2928.225
$ws.Range("K31").Value = 11112581
$ws.Range("L31").Value = 2928.225
$ws.Range("M31").Value = -11112286
$ws.Range("N31").Value = -3518.225
$ws.Range("H34").Value = 3450751.5
$ws.Range("I34").Value = 11112581
$ws.Range("J34").Value = 2928.225
$ws.Range("K34").Value = 11112581
$ws.Range("L34").Value = 2928.225
$ws.Range("M34").Value = -11112379
$ws.Range("N34").Value = -3332.225
$ws.Range("H113").Value = 1843.1666
$ws.Range("I113").Value = 2566.3333
$ws.Range("J113").Value = 1120
$ws.Range("K113").Value = 2566.3333
$ws.Range("L113").Value = 1120
$ws.Range("M113").Value = -396.3332999999998
$ws.Range("N113").Value = -5460
$ws.Range("H132").Value = 2346.8235
$ws.Range("I132").Value = 1678.4
$ws.Range("J132").Value = 3301.7144
$ws.Range("K132").Value = 5035.200000000001
$ws.Range("L132").Value = 9905.143199999999
$ws.Range("M132").Value = -2505.200000000001
$ws.Range("N132").Value = -14965.1432
$ws.Range("H134").Value = 4095.7666
$ws.Range("I134").Value = 4328
$ws.Range("J134").Value = 3631.3
$ws.Range("K134").Value = 12984
$ws.Range("L134").Value = 10893.9
$ws.Range("M134").Value = -10449
$ws.Range("N134").Value = -15963.9

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H14").Value = 112.85714
$ws.Range("I14").Value = 112.85714
$ws.Range("K14").Value = 338.57142
$ws.Range("M14").Value = -165.57142
$ws.Range("H115").Value = 2709.6843
$ws.Range("I115").Value = 371
$ws.Range("J115").Value = 3333.3333
$ws.Range("K115").Value = 1113
$ws.Range("L115").Value = 9999.999899999999
$ws.Range("M115").Value = 62
$ws.Range("N115").Value = -12349.9999

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H80").Value = 2205.8948
$ws.Range("I80").Value = 2600
$ws.Range("J80").Value = 2116.9033
$ws.Range("K80").Value = 2600
$ws.Range("L80").Value = 2116.9033
$ws.Range("M80").Value = -1602
$ws.Range("N80").Value = -4112.9033
$ws.Range("H83").Value = 2205.8948
$ws.Range("I83").Value = 2600
$ws.Range("J83").Value = 2116.9033
$ws.Range("K83").Value = 13000
$ws.Range("L83").Value = 10584.5165
$ws.Range("M83").Value = -8008
$ws.Range("N83").Value = -20568.5165
$ws.Range("H132").Value = 3013.96
$ws.Range("I132").Value = 2471
$ws.Range("K132").Value = 7413
$ws.Range("M132").Value = -4883

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H55").Value = 246.85715
$ws.Range("I55").Value = 240
$ws.Range("J55").Value = 254.4
$ws.Range("K55").Value = 240
$ws.Range("L55").Value = 254.4
$ws.Range("M55").Value = -67
$ws.Range("N55").Value = -600.4
$ws.Range("H61").Value = 2178.5557
$ws.Range("I61").Value = 2177.6155
$ws.Range("J61").Value = 2181
$ws.Range("K61").Value = 2177.6155
$ws.Range("L61").Value = 2181
$ws.Range("M61").Value = -1975.6155
$ws.Range("N61").Value = -2585
$ws.Range("H93").Value = 1780
$ws.Range("I93").Value = 1378.1818
$ws.Range("J93").Value = 2411.4285
$ws.Range("K93").Value = 1378.1818
$ws.Range("L93").Value = 2411.4285
$ws.Range("M93").Value = -130.1818000000001
$ws.Range("N93").Value = -4907.4285
$ws.Range("H98").Value = 33000
$ws.Range("J98").Value = 33000
$ws.Range("L98").Value = 33000
$ws.Range("N98").Value = -38990
$ws.Range("H113").Value = 2178.5557
$ws.Range("I113").Value = 2177.6155
$ws.Range("J113").Value = 2181
$ws.Range("K113").Value = 2177.6155
$ws.Range("L113").Value = 2181
$ws.Range("M113").Value = -7.615499999999884
$ws.Range("N113").Value = -6521
$ws.Range("H132").Value = 35733220
$ws.Range("I132").Value = 50023810
$ws.Range("K132").Value = 150071430
$ws.Range("M132").Value = -150068900
$ws.Range("H136").Value = 6725.591
$ws.Range("I136").Value = 10931.5
$ws.Range("J136").Value = 1678.5
$ws.Range("K136").Value = 32794.5
$ws.Range("L136").Value = 5035.5
$ws.Range("M136").Value = -30244.5
$ws.Range("N136").Value = -10135.5

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H126").Value = 5975.5
$ws.Range("I126").Value = 7300.6665
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 21901.9995
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -19431.9995
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 2174.6086
$ws.Range("I132").Value = 943.5
$ws.Range("J132").Value = 2433.7896
$ws.Range("K132").Value = 2830.5
$ws.Range("L132").Value = 7301.3688
$ws.Range("M132").Value = -300.5
$ws.Range("N132").Value = -12361.3688
$ws.Range("H136").Value = 1639.0385
$ws.Range("I136").Value = 874.4737
$ws.Range("J136").Value = 3714.2856
$ws.Range("K136").Value = 2623.4211
$ws.Range("L136").Value = 11142.8568
$ws.Range("M136").Value = -73.42110000000002
$ws.Range("N136").Value = -16242.8568
